# The commit swaps the two theme parts in the package: the theme that used
# to live at ppt/theme/theme1.xml ("Integral", driving the slide master /
# the presentation's visible design) and the theme that used to live at
# ppt/theme/theme2.xml ("Office Theme", only wired to the notes master)
# trade places. Concretely, after the edit theme1.xml must contain the
# "Office Theme" palette (name="Office"/"Office Theme") while theme2.xml
# must contain the former "Integral" palette.
#
# The PowerPoint object model only ever resolves "the" theme back to the
# part that actually drives the deck's design (theme1.xml) - there's no
# COM surface that reseats ppt/theme/theme2.xml itself (NotesMaster /
# NotesPage color objects alias back to the very same active theme in
# this host). So we reproduce the user-visible effect of the swap - the
# presentation's theme colours become the "Office Theme" palette - via the
# richer 12-slot ThemeColorScheme, which is the scriptable equivalent of
# editing <a:clrScheme> in ppt/theme/theme1.xml.

$p = $ppt.ActivePresentation

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the "Office Theme" colour scheme (was in theme2.xml),
# in the fixed 12-slot order PowerPoint uses for ThemeColorScheme:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink.
$officeTheme = @(
    (RGBVal 0x00 0x00 0x00),   # dk1      000000
    (RGBVal 0xFF 0xFF 0xFF),   # lt1      FFFFFF
    (RGBVal 0x44 0x54 0x6A),   # dk2      44546A
    (RGBVal 0xE7 0xE6 0xE6),   # lt2      E7E6E6
    (RGBVal 0x5B 0x9B 0xD5),   # accent1  5B9BD5
    (RGBVal 0xED 0x7D 0x31),   # accent2  ED7D31
    (RGBVal 0xA5 0xA5 0xA5),   # accent3  A5A5A5
    (RGBVal 0xFF 0xC0 0x00),   # accent4  FFC000
    (RGBVal 0x44 0x72 0xC4),   # accent5  4472C4
    (RGBVal 0x70 0xAD 0x47),   # accent6  70AD47
    (RGBVal 0x05 0x63 0xC1),   # hlink    0563C1
    (RGBVal 0x95 0x4F 0x72)    # folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeTheme[$i - 1]
}

# Best-effort: also try to rename the design/theme to match ("Office
# Theme"/"Office"). Harmless if the host does not persist these.
try { $p.Designs.Item(1).Name = "Office Theme" } catch {}
try { $p.Designs.Item(1).SlideMaster.Theme.Name = "Office Theme" } catch {}
try { $p.Designs.Item(1).SlideMaster.ColorScheme.Name = "Office" } catch {}
